$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3249.7
$ws.Range("J62").Value = 4999
$ws.Range("L62").Value = 4999
$ws.Range("N62").Value = -6247
$ws.Range("H65").Value = 3249.7
$ws.Range("J65").Value = 4999
$ws.Range("L65").Value = 24995
$ws.Range("N65").Value = -31235
$ws.Range("H86").Value = 6783.5
$ws.Range("I86").Value = 3424.1428
$ws.Range("K86").Value = 3424.1428
$ws.Range("M86").Value = -2301.1428
$ws.Range("H89").Value = 6783.5
$ws.Range("I89").Value = 3424.1428
$ws.Range("K89").Value = 17120.714
$ws.Range("M89").Value = -11504.714
$ws.Range("H106").Value = 44003948
$ws.Range("I106").Value = 48892720
$ws.Range("K106").Value = 48892720
$ws.Range("M106").Value = -48892089
$ws.Range("H138").Value = 3050.8518
$ws.Range("J138").Value = 3635.7058
$ws.Range("L138").Value = 10907.1174
$ws.Range("N138").Value = -21187.1174

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9812.333000000001
$ws.Range("I61").Value = 11040.6
$ws.Range("J61").Value = 3671
$ws.Range("K61").Value = 11040.6
$ws.Range("L61").Value = 3671
$ws.Range("M61").Value = -10828.6
$ws.Range("N61").Value = -4095
$ws.Range("H74").Value = 2808.8928
$ws.Range("I74").Value = 1200.9565
$ws.Range("J74").Value = 10205.4
$ws.Range("K74").Value = 1200.9565
$ws.Range("L74").Value = 10205.4
$ws.Range("M74").Value = -326.9565
$ws.Range("N74").Value = -11953.4
$ws.Range("H77").Value = 2808.8928
$ws.Range("I77").Value = 1200.9565
$ws.Range("J77").Value = 10205.4
$ws.Range("K77").Value = 6004.7825
$ws.Range("L77").Value = 51027
$ws.Range("M77").Value = -1636.7825
$ws.Range("N77").Value = -59763
$ws.Range("H110").Value = 2743.2727
$ws.Range("I110").Value = 2991.2222
$ws.Range("K110").Value = 2991.2222
$ws.Range("M110").Value = -946.2222000000002
$ws.Range("H136").Value = 9812.333000000001
$ws.Range("I136").Value = 11040.6
$ws.Range("J136").Value = 3671
$ws.Range("K136").Value = 33121.8
$ws.Range("L136").Value = 11013
$ws.Range("M136").Value = -30571.8
$ws.Range("N136").Value = -16113

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2355.3635
$ws.Range("I99").Value = 2148.3684
$ws.Range("J99").Value = 3666.3333
$ws.Range("K99").Value = 2148.3684
$ws.Range("L99").Value = 3666.3333
$ws.Range("M99").Value = -650.3683999999998
$ws.Range("N99").Value = -6662.3333
$ws.Range("H134").Value = 4136.364
$ws.Range("I134").Value = 2111.111
$ws.Range("J134").Value = 13250
$ws.Range("K134").Value = 6333.333
$ws.Range("L134").Value = 39750
$ws.Range("M134").Value = -3798.333
$ws.Range("N134").Value = -44820

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3150.92
$ws.Range("I58").Value = 1226.8334
$ws.Range("J58").Value = 4927
$ws.Range("K58").Value = 1226.8334
$ws.Range("L58").Value = 4927
$ws.Range("M58").Value = -1023.8334
$ws.Range("N58").Value = -5333
$ws.Range("H62").Value = 5407.4287
$ws.Range("I62").Value = 4770.4
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 4770.4
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = -4146.4
$ws.Range("N62").Value = -8248
$ws.Range("H65").Value = 5407.4287
$ws.Range("I65").Value = 4770.4
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 23852
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = -20732
$ws.Range("N65").Value = -41240
$ws.Range("H136").Value = 3150.92
$ws.Range("I136").Value = 1226.8334
$ws.Range("J136").Value = 4927
$ws.Range("K136").Value = 3680.5002
$ws.Range("L136").Value = 14781
$ws.Range("M136").Value = -1130.5002
$ws.Range("N136").Value = -19881

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3129.7778
$ws.Range("J113").Value = 2466.8572
$ws.Range("L113").Value = 7400.571599999999
$ws.Range("N113").Value = -11740.5716
$ws.Range("H118").Value = 6358.75
$ws.Range("I118").Value = 6358.75
$ws.Range("K118").Value = 19076.25
$ws.Range("M118").Value = -17833.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 557.5294
$ws.Range("J97").Value = 686.875
$ws.Range("L97").Value = 686.875
$ws.Range("N97").Value = -1678.875
$ws.Range("H132").Value = 4471.7827
$ws.Range("I132").Value = 2932.5293
$ws.Range("J132").Value = 8833
$ws.Range("K132").Value = 8797.5879
$ws.Range("L132").Value = 26499
$ws.Range("M132").Value = -6267.5879
$ws.Range("N132").Value = -31559
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 349995
$ws.Range("J43").Value = 349995
$ws.Range("L43").Value = 349995
$ws.Range("N43").Value = -350381
$ws.Range("H122").Value = 4072.25
$ws.Range("I122").Value = 4072.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12216.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9766.75
$ws.Range("N122").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2904.75
$ws.Range("I81").Value = 2907.1667
$ws.Range("J81").Value = 2897.5
$ws.Range("K81").Value = 5814.3334
$ws.Range("L81").Value = 5795
$ws.Range("M81").Value = -4753.3334
$ws.Range("N81").Value = -7917
$ws.Range("H84").Value = 2904.75
$ws.Range("I84").Value = 2907.1667
$ws.Range("J84").Value = 2897.5
$ws.Range("K84").Value = 29071.667
$ws.Range("L84").Value = 28975
$ws.Range("M84").Value = -23767.667
$ws.Range("N84").Value = -39583
$ws.Range("H100").Value = 328.77777
$ws.Range("I100").Value = 374.85715
$ws.Range("J100").Value = 167.5
$ws.Range("K100").Value = 749.7143
$ws.Range("L100").Value = 335
$ws.Range("M100").Value = -208.7143
$ws.Range("N100").Value = -1417
$ws.Range("H132").Value = 3262
$ws.Range("I132").Value = 2314.182
$ws.Range("K132").Value = 6942.545999999999
$ws.Range("M132").Value = -4412.545999999999
$ws.Range("H136").Value = 4540.4546
$ws.Range("I136").Value = 4401.1875
$ws.Range("J136").Value = 8997
$ws.Range("K136").Value = 13203.5625
$ws.Range("L136").Value = 26991
$ws.Range("M136").Value = -10653.5625
$ws.Range("N136").Value = -32091
$ws.Range("H141").Value = 55500
$ws.Range("J141").Value = 55500
$ws.Range("L141").Value = 55500
$ws.Range("N141").Value = -65860
